# Generate Report for Handback
# Adds a new handback record (908ad72e-1216-4cf0-8b8f-d81bd173824d.md) as
# row 4 to the "Overview", "zh-cn" and "de-de" worksheets/tables.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$hyperFontColor = 15570276   # BGR packed value of RGB(0x64,0x95,0xED) == FF6495ED
$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Format-HyperlinkCell($cell) {
    $cell.Font.Underline = 2
    $cell.Font.Color = $hyperFontColor
    $cell.Font.Name = "Calibri"
}

# ---------------------------------------------------------------------------
# Sheet "Overview" - row 4
# ---------------------------------------------------------------------------
$wsOverview.Cells.Item(4, 1).Value = "908ad72e-1216-4cf0-8b8f-d81bd173824d.md"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item(4, 2),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/372aa4df753b765d5f043c4869e8221439b26ea/e2e/908ad72e-1216-4cf0-8b8f-d81bd173824d.md",
    "",
    "",
    "e2e\908ad72e-1216-4cf0-8b8f-d81bd173824d.md"
)
Format-HyperlinkCell $wsOverview.Cells.Item(4, 2)

$wsOverview.Cells.Item(4, 3).Value = ".md"
$wsOverview.Cells.Item(4, 5).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(4, 6).Value = "Handed back: in sync with en-US"

$wsOverview.Cells.Item(4, 7).Value = "2016-08-28 04:44:07"
$wsOverview.Cells.Item(4, 7).NumberFormat = $dateFormat

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" - row 4
# ---------------------------------------------------------------------------
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item(4, 1),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/372aa4df753b765d5f043c4869e8221439b26ea/e2e/908ad72e-1216-4cf0-8b8f-d81bd173824d.md",
    "",
    "",
    "908ad72e-1216-4cf0-8b8f-d81bd173824d.md"
)
Format-HyperlinkCell $wsZhCn.Cells.Item(4, 1)

$wsZhCn.Cells.Item(4, 2).Value = ".md"
$wsZhCn.Cells.Item(4, 3).Value = "Handed back: in sync with en-US"
$wsZhCn.Cells.Item(4, 4).Value = "e2e"
$wsZhCn.Cells.Item(4, 5).Value = "ht"
$wsZhCn.Cells.Item(4, 6).Value = "True"
$wsZhCn.Cells.Item(4, 7).Value = "908ad72e-1216-4cf0-8b8f-d81bd173824d.bfdf3c9572a38d86c814029a5685277c177c5877.zh-cn.xlf"

$wsZhCn.Cells.Item(4, 8).Value = "2016-08-28 04:43:59"
$wsZhCn.Cells.Item(4, 8).NumberFormat = $dateFormat

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item(4, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/775a8d3dae959451cd1417931039bba16e9ed46/e2e/908ad72e-1216-4cf0-8b8f-d81bd173824d.md",
    "",
    "",
    "908ad72e-1216-4cf0-8b8f-d81bd173824d.md"
)
Format-HyperlinkCell $wsZhCn.Cells.Item(4, 9)

$wsZhCn.Cells.Item(4, 10).Value = "908ad72e-1216-4cf0-8b8f-d81bd173824d.bfdf3c9572a38d86c814029a5685277c177c5877.zh-cn.xlf"

$wsZhCn.Cells.Item(4, 11).Value = "2016-08-28 04:44:27"
$wsZhCn.Cells.Item(4, 11).NumberFormat = $dateFormat

$wsZhCn.Cells.Item(4, 12).Value = ""
$wsZhCn.Cells.Item(4, 13).Value = "True"
$wsZhCn.Cells.Item(4, 14).Value = ""
$wsZhCn.Cells.Item(4, 15).Value = "False"
$wsZhCn.Cells.Item(4, 16).Value = ""

$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------------
# Sheet "de-de" - row 4
# ---------------------------------------------------------------------------
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item(4, 1),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/372aa4df753b765d5f043c4869e8221439b26ea/e2e/908ad72e-1216-4cf0-8b8f-d81bd173824d.md",
    "",
    "",
    "908ad72e-1216-4cf0-8b8f-d81bd173824d.md"
)
Format-HyperlinkCell $wsDeDe.Cells.Item(4, 1)

$wsDeDe.Cells.Item(4, 2).Value = ".md"
$wsDeDe.Cells.Item(4, 3).Value = "Handed back: in sync with en-US"
$wsDeDe.Cells.Item(4, 4).Value = "e2e"
$wsDeDe.Cells.Item(4, 5).Value = "ht"
$wsDeDe.Cells.Item(4, 6).Value = "True"
$wsDeDe.Cells.Item(4, 7).Value = "908ad72e-1216-4cf0-8b8f-d81bd173824d.bfdf3c9572a38d86c814029a5685277c177c5877.de-de.xlf"

$wsDeDe.Cells.Item(4, 8).Value = "2016-08-28 04:44:07"
$wsDeDe.Cells.Item(4, 8).NumberFormat = $dateFormat

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item(4, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0bbed0ff057c45ff0b2e1cd8a50f0fa10c43983/e2e/908ad72e-1216-4cf0-8b8f-d81bd173824d.md",
    "",
    "",
    "908ad72e-1216-4cf0-8b8f-d81bd173824d.md"
)
Format-HyperlinkCell $wsDeDe.Cells.Item(4, 9)

$wsDeDe.Cells.Item(4, 10).Value = "908ad72e-1216-4cf0-8b8f-d81bd173824d.bfdf3c9572a38d86c814029a5685277c177c5877.de-de.xlf"

$wsDeDe.Cells.Item(4, 11).Value = "2016-08-28 04:44:33"
$wsDeDe.Cells.Item(4, 11).NumberFormat = $dateFormat

$wsDeDe.Cells.Item(4, 12).Value = ""
$wsDeDe.Cells.Item(4, 13).Value = "True"
$wsDeDe.Cells.Item(4, 14).Value = ""
$wsDeDe.Cells.Item(4, 15).Value = "False"
$wsDeDe.Cells.Item(4, 16).Value = ""

$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:P4"))
